$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date-serial value (YYYY-MM-DD formatted) for each
# data row (rows 2-28). This automatic update bumps that date by one day
# (45444 -> 45445, i.e. 2024-06-01 -> 2024-06-02) for every row.
for ($row = 2; $row -le 28; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45444) {
        $cell.Value2 = 45445
    }
}
